# "Tried to implement Penalty Reward System (unfinished)"
# - Shifts the weekly forecast dates forward by one week (each week's
#   Week_Start_Date becomes the following week's date) and updates the
#   MyForecast values on the "Forecast Comparison" sheet.
# - Recomputes the dependent summary statistics on the "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Column B on "Forecast Comparison" holds plain-text dates (e.g. "2025-01-12")
# rather than real Excel dates - force text formatting before writing so the
# values don't get auto-converted into date serial numbers.
$ws1.Range("B2:B17").NumberFormat = "@"

# Week_Start_Date (col B) + MyForecast (col D) per row.
$ws1.Range("B2").Value  = "2025-01-12"
$ws1.Range("D2").Value  = 187

$ws1.Range("B3").Value  = "2025-01-19"
$ws1.Range("D3").Value  = 203

$ws1.Range("B4").Value  = "2025-01-26"
$ws1.Range("D4").Value  = 199

$ws1.Range("B5").Value  = "2025-02-02"
$ws1.Range("D5").Value  = 182

$ws1.Range("B6").Value  = "2025-02-09"
$ws1.Range("D6").Value  = 167

$ws1.Range("B7").Value  = "2025-02-16"
$ws1.Range("D7").Value  = 169

$ws1.Range("B8").Value  = "2025-02-23"
$ws1.Range("D8").Value  = 190

$ws1.Range("B9").Value  = "2025-03-02"
$ws1.Range("D9").Value  = 143

$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("D10").Value = 138

$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("D11").Value = 134

$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("D12").Value = 187

$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("D13").Value = 180

$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("D14").Value = 128

$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("D15").Value = 125

$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("D16").Value = 120

$ws1.Range("B17").Value = "2025-04-27"
$ws1.Range("D17").Value = 123

# Summary sheet values are stored as plain text too - force text formatting
# so numeric-looking strings ("303", "133", ...) aren't coerced to numbers.
$ws2.Range("B2").NumberFormat  = "@"
$ws2.Range("B4").NumberFormat  = "@"
$ws2.Range("B6").NumberFormat  = "@"
$ws2.Range("B8").NumberFormat  = "@"
$ws2.Range("B9").NumberFormat  = "@"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B15").NumberFormat = "@"

$ws2.Range("B2").Value  = "2022-12-25 to 2025-01-05"   # Historical Range
$ws2.Range("B4").Value  = "303"                         # Max Sales
$ws2.Range("B6").Value  = "133"                         # Median Sales
$ws2.Range("B8").Value  = "14098 units"                 # Total Historical Sales
$ws2.Range("B9").Value  = "2574"                        # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "1440"                        # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "771"                         # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "203"                         # Max Forecast
$ws2.Range("B13").Value = "2025-01-19"                  # Max Forecast Week
$ws2.Range("B14").Value = "120"                         # Min Forecast
$ws2.Range("B15").Value = "2025-04-20"                  # Min Forecast Week
